# Adds a "midi" column to the primary cue/monitor/visual block and a
# "visual" column to the mc1: group, then records the image filenames
# used for the basecamp stage (content-uri config for initial images).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before L and label it "midi"
# (shifts the old L:Z block one column to the right, to M:AA)
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("L1").Value = "midi"

# Insert a new column before T (post-shift coordinates) and label it "visual"
# (shifts the old T:AA block one column to the right, to U:AB)
$ws.Range("T1").EntireColumn.Insert()
$ws.Range("T1").Value = "visual"

# Rename the basecamp image filename to match the on-disk asset casing
$ws.Range("N2").Value = "baseCamp.jpg"

# Record the new "ending" image filename in the newly inserted visual column
$ws.Range("T2").Value = "ending.jpg"

# Update the view to match: scrolled right with L1 selected/active
$ws.Range("L1").Select()
